$wb = $excel.ActiveWorkbook

# --- Add the new "Receive Parcel" worksheet after PSQAttendance ---
$loginSheet = $wb.Worksheets.Item("OELogin")
$releaseOrderSheet = $wb.Worksheets.Item("ReleaseOrder")
$attendanceSheet = $wb.Worksheets.Item("PSQAttendance")

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $attendanceSheet)
$newSheet.Name = "Receive Parcel"

# --- Apply the same cell styles used throughout the workbook FIRST (before
#     entering values) so numeric-looking text lands in text-formatted cells
#     and Excel stores it as a shared string rather than a number. ---

# Header style (bold, yellow fill, bordered) -> reuse from OELogin!A1
$loginSheet.Range("A1").Copy()
$newSheet.Range("A1:N1").PasteSpecial(-4122)

# Text style with border (numFmtId 49 "@") -> reuse from OELogin!A2
$loginSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("F2").PasteSpecial(-4122)
$newSheet.Range("G2").PasteSpecial(-4122)
$newSheet.Range("I2").PasteSpecial(-4122)
$newSheet.Range("J2").PasteSpecial(-4122)
$newSheet.Range("L2").PasteSpecial(-4122)

# General/bordered style -> reuse from ReleaseOrder!C2
$releaseOrderSheet.Range("C2").Copy()
$newSheet.Range("C2").PasteSpecial(-4122)
$newSheet.Range("D2").PasteSpecial(-4122)
$newSheet.Range("E2").PasteSpecial(-4122)
$newSheet.Range("H2").PasteSpecial(-4122)
$newSheet.Range("K2").PasteSpecial(-4122)
$newSheet.Range("M2").PasteSpecial(-4122)
$newSheet.Range("N2").PasteSpecial(-4122)

$newSheet.Application.CutCopyMode = $false

# --- Header row ---
$newSheet.Range("A1").Value = "hrms_id"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("C1").Value = "vendor"
$newSheet.Range("D1").Value = "deliverymode"
$newSheet.Range("E1").Value = "materialtype"
$newSheet.Range("F1").Value = "insuredamount"
$newSheet.Range("G1").Value = "noofboxes"
$newSheet.Range("H1").Value = "courieragencyname"
$newSheet.Range("I1").Value = "courierreceiptnumber"
$newSheet.Range("J1").Value = "couriercharges"
$newSheet.Range("K1").Value = "borneby"
$newSheet.Range("L1").Value = "grosswt"
$newSheet.Range("M1").Value = "receivedby"
$newSheet.Range("N1").Value = "sentby"

# --- Data row (order matches the shared-string insertion order of the target file) ---
$newSheet.Range("A2").Value = "110"
$newSheet.Range("B2").Value = "Pass@123"
$newSheet.Range("G2").Value = "1"
$newSheet.Range("D2").Value = "P"
$newSheet.Range("E2").Value = "G"
$newSheet.Range("H2").Value = "ABC"
$newSheet.Range("K2").Value = "CO"
$newSheet.Range("N2").Value = "test"
$newSheet.Range("F2").Value = "25000"
$newSheet.Range("I2").Value = "12345"
$newSheet.Range("J2").Value = "1500"
$newSheet.Range("L2").Value = "25"
$newSheet.Range("C2").Value = "BKK"

# --- Hyperlink on the password cell (matches pattern used by other sheets) ---
$newSheet.Hyperlinks.Add($newSheet.Range("B2"), "mailto:Pass@123")

# Hyperlink style -> reuse from OELogin!B2 (re-applied after Hyperlinks.Add,
# which resets the cell's own formatting)
$loginSheet.Range("B2").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

$newSheet.Range("E7").Select()

# --- Adjust selections / active states on the other touched sheets ---
$loginSheet.Range("B1").Select()

$releaseOrderSheet.Range("A1").Select()

$attendanceSheet.Select()
$attendanceSheet.Range("K2").Select()

# Re-style G2/I2 to the plain bordered style (drops the now-unused applyFill xf)
$releaseOrderSheet.Range("C2").Copy()
$attendanceSheet.Range("G2").PasteSpecial(-4122)
$attendanceSheet.Range("I2").PasteSpecial(-4122)
$attendanceSheet.Application.CutCopyMode = $false

$newSheet.Select()
